$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.835.77'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '2.315.86'
$ws.Range('E3').Value = '  +1.59%  '
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '114.19'
$ws.Range('E5').Value = '  +19.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '270.80'
$ws.Range('E6').Value = '  +1.38%  '
$ws.Range('E7').Value = '  +0.64%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.623'
$ws.Range('E9').Value = '  +2.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '47.47'
$ws.Range('E10').Value = '  +7.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0943'
$ws.Range('E11').Value = '  +0.84%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.90'
$ws.Range('E12').Value = '  +15.51%  '
$ws.Range('E13').Value = '  +1.76%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.76'
$ws.Range('E14').Value = '  +3.93%  '
$ws.Range('D15').Value = '2.662.69'
$ws.Range('E15').Value = '  +1.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.862'
$ws.Range('E16').Value = '  +1.92%  '
$ws.Range('D17').Value = '2.320.97'
$ws.Range('E17').Value = '  +1.72%  '
$ws.Range('D18').Value = '43.828.82'
$ws.Range('E18').Value = '  +0.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000110'
$ws.Range('E19').Value = '  +3.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.68'
$ws.Range('E20').Value = '  +8.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.82'
$ws.Range('E21').Value = '  +1.24%  '
$ws.Range('E22').Value = '  +5.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '234.51'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.91'
$ws.Range('E24').Value = '  +16.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.45'
$ws.Range('E25').Value = '  +5.88%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').Value = '  +0.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '42.84'
$ws.Range('E28').Value = '  +10.54%  '
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '177.87'
$ws.Range('E31').Value = '  +0.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.97'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0934'
$ws.Range('E33').Value = '  +5.59%  '
$ws.Range('E34').Value = '  +4.40%  '
$ws.Range('E35').Value = '  +0.77%  '
$ws.Range('E36').Value = '  +6.61%  '
$ws.Range('E37').Value = '  +3.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.98'
$ws.Range('E38').Value = '  +21.48%  '
$ws.Range('E39').Value = '  +0.27%  '
$ws.Range('E40').Value = '  +3.81%  '
$ws.Range('E41').Value = '  +1.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '70.35'
$ws.Range('E42').Value = '  +12.47%  '
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('B44').Value = 'THORChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.98'
$ws.Range('E44').Value = '  +14.69%  '
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.40'
$ws.Range('E45').Value = '  +3.41%  '
$ws.Range('B46').Value = 'Celestia'
$ws.Range('C46').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.67'
$ws.Range('E46').Value = '  +6.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.82'
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('E48').Value = '  -0.77%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.475'
$ws.Range('E49').Value = '  +12.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '100.48'
$ws.Range('E50').Value = '  +2.15%  '
$ws.Range('E51').Value = '  +3.56%  '
